$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Value)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "26.327.91"
Set-TextValue $ws.Range("E2") "  +0.31%  "
Set-TextValue $ws.Range("D3") "1.667.80"
Set-TextValue $ws.Range("E3") "  +0.57%  "
Set-TextValue $ws.Range("D4") "1.008"
Set-TextValue $ws.Range("E4") "  -0.02%  "
Set-TextValue $ws.Range("D5") "220.61"
Set-TextValue $ws.Range("E5") "  +1.03%  "
Set-TextValue $ws.Range("D6") "0.5315"
Set-TextValue $ws.Range("E6") "  -0.07%  "
Set-TextValue $ws.Range("E7") "  -0.02%  "
Set-TextValue $ws.Range("D8") "0.2649"
Set-TextValue $ws.Range("E8") "  +0.78%  "
Set-TextValue $ws.Range("E9") "  +0.25%  "
Set-TextValue $ws.Range("D10") "20.94"
Set-TextValue $ws.Range("E10") "  +2.14%  "
Set-TextValue $ws.Range("D11") "0.07839"
Set-TextValue $ws.Range("E11") "  -0.13%  "
Set-TextValue $ws.Range("D12") "4.536"
Set-TextValue $ws.Range("E12") "  -0.03%  "
Set-TextValue $ws.Range("D13") "1.671.39"
Set-TextValue $ws.Range("E13") "  +0.70%  "
Set-TextValue $ws.Range("D14") "1.897.75"
Set-TextValue $ws.Range("E14") "  +0.65%  "
Set-TextValue $ws.Range("D15") "0.5614"
Set-TextValue $ws.Range("E15") "  +1.88%  "
Set-TextValue $ws.Range("D16") "0.0₅8147"
Set-TextValue $ws.Range("E16") "  -0.41%  "
Set-TextValue $ws.Range("D17") "65.87"
Set-TextValue $ws.Range("E17") "  +0.50%  "
Set-TextValue $ws.Range("D18") "26.328.75"
Set-TextValue $ws.Range("E18") "  +0.42%  "
Set-TextValue $ws.Range("D19") "1.008"
Set-TextValue $ws.Range("E19") "  -0.03%  "
Set-TextValue $ws.Range("D20") "4.722"
Set-TextValue $ws.Range("E20") "  +2.20%  "
Set-TextValue $ws.Range("D21") "197.96"
Set-TextValue $ws.Range("E21") "  +3.15%  "
Set-TextValue $ws.Range("E22") "  +1.59%  "
Set-TextValue $ws.Range("D23") "6.058"
Set-TextValue $ws.Range("E23") "  +0.46%  "
Set-TextValue $ws.Range("E24") "  -0.01%  "
Set-TextValue $ws.Range("D25") "146.62"
Set-TextValue $ws.Range("E25") "  +1.99%  "
Set-TextValue $ws.Range("E26") "  -0.34%  "
Set-TextValue $ws.Range("D27") "7.256"
Set-TextValue $ws.Range("E27") "  +0.36%  "
Set-TextValue $ws.Range("E28") "  +0.65%  "
Set-TextValue $ws.Range("D29") "1.510"
Set-TextValue $ws.Range("E29") "  +2.68%  "
Set-TextValue $ws.Range("D30") "0.05900"
Set-TextValue $ws.Range("E30") "  +1.93%  "
Set-TextValue $ws.Range("D31") "1.284"
Set-TextValue $ws.Range("E31") "  +0.58%  "
Set-TextValue $ws.Range("D32") "3.553"
Set-TextValue $ws.Range("E32") "  -0.35%  "
Set-TextValue $ws.Range("D33") "3.329"
Set-TextValue $ws.Range("E33") "  +1.39%  "
Set-TextValue $ws.Range("D34") "1.605"
Set-TextValue $ws.Range("E34") "  +0.29%  "
Set-TextValue $ws.Range("B35") "MXToken"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D35") "2.835"
Set-TextValue $ws.Range("E35") "  +0.65%  "
Set-TextValue $ws.Range("B36") "ARBITRUM"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D36") "0.9622"
Set-TextValue $ws.Range("E36") "  +0.89%  "
Set-TextValue $ws.Range("D37") "2.434"
Set-TextValue $ws.Range("E37") "  +0.27%  "
Set-TextValue $ws.Range("D38") "0.5818"
Set-TextValue $ws.Range("E38") "  +0.53%  "
Set-TextValue $ws.Range("D39") "0.01619"
Set-TextValue $ws.Range("E39") "  +0.91%  "
Set-TextValue $ws.Range("D40") "5.957"
Set-TextValue $ws.Range("E40") "  +2.34%  "
Set-TextValue $ws.Range("D41") "1.074.88"
Set-TextValue $ws.Range("E41") "  +3.00%  "
Set-TextValue $ws.Range("D42") "0.8573"
Set-TextValue $ws.Range("E42") "  +0.71%  "
Set-TextValue $ws.Range("D43") "1.008"
Set-TextValue $ws.Range("E43") "  -0.05%  "
Set-TextValue $ws.Range("D44") "102.77"
Set-TextValue $ws.Range("E44") "  -1.82%  "
Set-TextValue $ws.Range("D45") "1.808.13"
Set-TextValue $ws.Range("E45") "  +0.52%  "
Set-TextValue $ws.Range("D46") "58.57"
Set-TextValue $ws.Range("E46") "  +2.91%  "
Set-TextValue $ws.Range("D47") "1.013"
Set-TextValue $ws.Range("E47") "  +0.47%  "
Set-TextValue $ws.Range("B48") "Mantle"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D48") "0.4407"
Set-TextValue $ws.Range("E48") "  +0.83%  "
Set-TextValue $ws.Range("B49") "EnergySwap"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "8.076"
Set-TextValue $ws.Range("E49") "  +2.33%  "
Set-TextValue $ws.Range("B50") "BabyDogeCoin"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D50") "0.0₈102"
Set-TextValue $ws.Range("E50") "  -2.81%  "
Set-TextValue $ws.Range("D51") "0.05149"
Set-TextValue $ws.Range("E51") "  -0.19%  "
